# Roll back the recent changes to the "InvalidLogin" worksheet:
#  - restore the original A1:B5 demo credential values
#  - remove the extra data that had been added in R21:S25
#  - restore the selection to A1:B5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvalidLogin")

# Restore original values for the A1:B5 block
$ws.Range("A1").Value = "admin123"
$ws.Range("B1").Value = "manager123"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager123"
$ws.Range("A3").Value = "admin123"
$ws.Range("B3").Value = "manager"
$ws.Range("A4").Value = "admin"
$ws.Range("B4").ClearContents()
$ws.Range("A5").ClearContents()
$ws.Range("B5").Value = "manager"

# Remove the extra rows 21:25 (columns R:S) that had been added
$ws.Range("R21:S25").ClearContents()

# Restore the original selection covering A1:B5
$ws.Range("A1:B5").Select()
